# "Added cards to every user" - fill in the Card # column (E) for every
# user row on the UserList sheet that didn't have one yet, then leave the
# selection where the last edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserList")

$ws.Range("E3").Value = 100000001
$ws.Range("E4").Value = 100000002
$ws.Range("E5").Value = 100000003
$ws.Range("E6").Value = 100000004
$ws.Range("E7").Value = 100000005
$ws.Range("E8").Value = 100000006
$ws.Range("E9").Value = 100000007

$ws.Range("E9").Select()
